$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; D = "298.20"; E = "-1.89%" },
    @{ Row = 3; D = "31.35"; E = "-1.55%" },
    @{ Row = 4; D = "5.121"; E = "-2.56%" },
    @{ Row = 5; D = "0.07331"; E = "-2.62%" },
    @{ Row = 6; D = "7.751"; E = "-0.97%" },
    @{ Row = 7; D = "1.747"; E = "17.94%" },
    @{ Row = 8; D = "3.727"; E = "-0.83%" },
    @{ Row = 9; D = "0.9241"; E = "1.08%" },
    @{ Row = 10; D = "0.1665"; E = "-1.74%" },
    @{ Row = 11; D = "0.06907"; E = "-8.55%" },
    @{ Row = 12; D = "0.07940"; E = "-0.19%" },
    @{ Row = 13; D = "0.02999"; E = "0.19%" },
    @{ Row = 14; D = "0.09911"; E = "0.20%" },
    @{ Row = 15; D = "0.001493"; E = "-0.35%" },
    @{ Row = 16; D = "0.006252"; E = "-0.22%" },
    @{ Row = 17; D = "3.455"; E = "-1.15%" },
    @{ Row = 18; D = "2.221"; E = "-0.48%" },
    @{ Row = 19; E = "-2.55%" },
    @{ Row = 20; D = "0.1316"; E = "-2.11%" },
    @{ Row = 21; D = "4.549"; E = "1.71%" },
    @{ Row = 22; D = "0.04634"; E = "1.68%" },
    @{ Row = 24; D = "0.001220"; E = "0.36%" },
    @{ Row = 25; D = "0.004740"; E = "6.69%" },
    @{ Row = 26; D = "0.0001297"; E = "-7.17%" },
    @{ Row = 27; D = "0.0001872"; E = "7.74%" },
    @{ Row = 39; D = "0.01721"; E = "4.34%" },
    @{ Row = 40; D = "0.04454"; E = "-0.97%" },
    @{ Row = 41; D = "0.007128"; E = "-1.24%" },
    @{ Row = 42; D = "0.1331"; E = "-1.19%" },
    @{ Row = 43; D = "0.002204"; E = "-1.81%" },
    @{ Row = 44; D = "0.01080"; E = "-15.84%" },
    @{ Row = 45; D = "0.00006158"; E = "-0.85%" },
    @{ Row = 46; E = "-21.27%" },
    @{ Row = 47; D = "0.7388"; E = "4.18%" }
)

foreach ($chg in $changes) {
    $row = $chg.Row
    if ($chg.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $chg.D
        $cell.Style = "Normal"
    }
    if ($chg.ContainsKey("E")) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $chg.E
        $cell.Style = "Normal"
    }
}
